# Config/Datas/ModelData.xlsx
# "StateMachine 绑定 Animation 名称 && Animation 名称进行 Mix 匹配"
#
# The ModelInfo sheet's "##" row documents field names/paths for the data
# rows below it. The field previously named "AnimationMix" (with Chinese
# description "动画混合路径") is renamed to "Animation" ("动画路径").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 = field name header, D3 = Chinese description of the field
$ws.Range("D1").Value = "Animation"
$ws.Range("D3").Value = "动画路径"

# Workbook default ("Normal") font switched from Calibri to 宋体 (SimSun)
$wb.Styles.Item("Normal").Font.Name = "宋体"
